# Apply updated crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.003.55'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '2.472.39'
$ws.Range("E3").Value = '  +2.10%  '
$ws.Range("E4").Value = '  -0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").Value = '2.471.67'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.71%  '
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("E16").Value = '  +1.95%  '
$ws.Range("D17").Value = '63.080.02'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '2.470.43'
$ws.Range("E18").Value = '  +0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.13%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '664.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.48%  '
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").Value = '2.590.02'
$ws.Range("E29").Value = '  +0.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("E33").Value = '  +0.50%  '
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("E35").Value = '  +2.87%  '
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("E37").Value = '  +0.49%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.371'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.41'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '151.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '0.0₆0306'
$ws.Range("E45").Value = '  +6.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +26.66%  '
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("E51").Value = '  -0.90%  '
